# Add a new "2021" data column (P) to the right of the existing "2020"
# column (O), mirroring each row's existing formatting, then leave the
# selection on Q4 (the cell to the right of the new column), matching
# the author's last on-screen selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-FormattedValue($sourceRef, $targetRef, $value) {
    $ws.Range($sourceRef).Copy()
    $ws.Range($targetRef).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    if ($null -ne $value) {
        $ws.Range($targetRef).Value = $value
    }
}

# Row 3: bottom border row under the header, no value (style like C3)
Set-FormattedValue "C3" "P3" $null

# Row 4: year header, style like O4
Set-FormattedValue "O4" "P4" 2021

# Row 5: "Total" row, numeric, style like O8
Set-FormattedValue "O8" "P5" 9038

# Row 6: section header row, no value, style like O6
Set-FormattedValue "O6" "P6" $null

# Rows 7-8: numeric rows, style like O9
Set-FormattedValue "O9" "P7" 8587
Set-FormattedValue "O9" "P8" 451

# Row 9: section header row, no value, style like O9
Set-FormattedValue "O9" "P9" $null

# Rows 10-24: "…" placeholder (no 2021 breakdown data yet), style like each O cell
10..24 | ForEach-Object {
    Set-FormattedValue "O$_" "P$_" "…"
}

# Row 25: bottom (bordered) row, "…" placeholder, style like O25
Set-FormattedValue "O25" "P25" "…"

# Leave the selection where the author left it when saving
$ws.Range("Q4").Select()
